$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.636.66'
$ws.Range("E2").Value = '  -1.04%  '

$ws.Range("D3").Value = '3.307.31'
$ws.Range("E3").Value = '  -0.05%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.25%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.09'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.19%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '183.96'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.59%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").Value = '3.301.04'
$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.571'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.65%  '

$ws.Range("E10").Value = '  -4.23%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.572'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.01%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '46.64'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.40%  '

$ws.Range("E13").Value = '  -2.33%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '634.54'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.36%  '

$ws.Range("D15").Value = '3.836.36'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.45'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.71%  '

$ws.Range("D17").Value = '65.793.99'
$ws.Range("E17").Value = '  -0.79%  '

$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.117'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.11%  '

$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.86'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.26%  '

$ws.Range("D20").Value = '3.304.67'
$ws.Range("E20").Value = '  -0.18%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.00'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.86%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.890'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.48%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.71'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.20%  '

$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.01'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.44%  '

$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '100.24'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.85%  '

$ws.Range("E26").Value = '  -0.39%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.72'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.33%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.40'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.72%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '30.85'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.68%  '

$ws.Range("E30").Value = '  -3.87%  '

$ws.Range("E31").Value = '  -2.14%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '594.53'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.62%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.74'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -8.91%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '10.94'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.64%  '

$ws.Range("D35").Value = '3.854.61'
$ws.Range("E35").Value = '  +3.28%  '

$ws.Range("E36").Value = '  -1.04%  '

$ws.Range("E37").Value = '  +0.14%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '55.69'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.22%  '

$ws.Range("D39").Value = '0.0₃0699'
$ws.Range("E39").Value = '  -4.81%  '

$ws.Range("E40").Value = '  -2.89%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.41'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.73%  '

$ws.Range("E42").Value = '  -5.18%  '

$ws.Range("E43").Value = '  -5.42%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.60'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.91%  '

$ws.Range("E45").Value = '  -2.59%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0408'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.93%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.03'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -10.46%  '

$ws.Range("E48").Value = '  -2.09%  '

$ws.Range("E49").Value = '  +0.52%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.52'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.81%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '130.70'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.53%  '
